$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column cells are stored as text in the workbook (t="inlineStr").
# Force every Price cell we touch to keep a Text number format before writing
# its new value so Excel does not silently reinterpret it as a floating point
# number (which would lose trailing zeros / exact decimal text, e.g. "1.00").
$dCells = @("D2","D3","D5","D6","D7","D8","D9","D10","D11","D12","D14","D15","D16","D17","D19","D20","D21","D22","D23","D24","D25","D26","D27","D28","D37","D38","D39","D40","D41","D43","D45","D47","D48","D50")
foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '57.007.36'
$ws.Range("E2").Value = '  +6.33%  '
$ws.Range("D3").Value = '3.236.98'
$ws.Range("E3").Value = '  +2.90%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '394.96'
$ws.Range("E5").Value = '  -0.49%  '
$ws.Range("D6").Value = '107.32'
$ws.Range("E6").Value = '  -2.35%  '
$ws.Range("D7").Value = '0.573'
$ws.Range("E7").Value = '  +4.75%  '
$ws.Range("D8").Value = '3.233.40'
$ws.Range("E8").Value = '  +2.90%  '
$ws.Range("D9").Value = '1.00'
$ws.Range("E9").Value = '  +0.02%  '
$ws.Range("D10").Value = '0.617'
$ws.Range("E10").Value = '  +0.92%  '
$ws.Range("D11").Value = '39.07'
$ws.Range("D12").Value = '0.0958'
$ws.Range("E12").Value = '  +9.83%  '
$ws.Range("D14").Value = '3.746.54'
$ws.Range("E14").Value = '  +2.56%  '
$ws.Range("D15").Value = '8.21'
$ws.Range("E15").Value = '  +2.34%  '
$ws.Range("D16").Value = '19.14'
$ws.Range("E16").Value = '  +0.32%  '
$ws.Range("D17").Value = '3.236.61'
$ws.Range("E17").Value = '  +2.62%  '
$ws.Range("E18").Value = '  -2.29%  '
$ws.Range("D19").Value = '10.87'
$ws.Range("E19").Value = '  +2.95%  '
$ws.Range("D20").Value = '56.839.24'
$ws.Range("E20").Value = '  +6.25%  '
$ws.Range("D21").Value = '3.34'
$ws.Range("E21").Value = '  +1.72%  '
$ws.Range("D22").Value = '0.0000106'
$ws.Range("E22").Value = '  +8.78%  '
$ws.Range("D23").Value = '12.96'
$ws.Range("E23").Value = '  +0.90%  '
$ws.Range("D24").Value = '296.75'
$ws.Range("E24").Value = '  +9.46%  '
$ws.Range("D25").Value = '74.01'
$ws.Range("E25").Value = '  +4.38%  '
$ws.Range("D26").Value = '3.17'
$ws.Range("E26").Value = '  -2.52%  '
$ws.Range("D27").Value = '27.83'
$ws.Range("E27").Value = '  +0.99%  '
$ws.Range("D28").Value = '7.69'
$ws.Range("E28").Value = '  -4.70%  '
$ws.Range("E29").Value = '  -1.19%  '
$ws.Range("E30").Value = '  -0.31%  '
$ws.Range("E32").Value = '  +3.77%  '
$ws.Range("E33").Value = '  -1.34%  '
$ws.Range("E34").Value = '  +2.14%  '
$ws.Range("E35").Value = '  -2.51%  '
$ws.Range("E36").Value = '  +1.64%  '
$ws.Range("D37").Value = '51.71'
$ws.Range("E37").Value = '  +2.41%  '
$ws.Range("D38").Value = '3.53'
$ws.Range("E38").Value = '  -3.41%  '
$ws.Range("D39").Value = '0.999'
$ws.Range("E39").Value = '  +0.00%  '
$ws.Range("D40").Value = '2.96'
$ws.Range("E40").Value = '  +5.87%  '
$ws.Range("D41").Value = '134.52'
$ws.Range("E41").Value = '  +3.20%  '
$ws.Range("E42").Value = '  +2.26%  '
$ws.Range("D43").Value = '3.97'
$ws.Range("E43").Value = '  -3.44%  '
$ws.Range("E44").Value = '  -0.28%  '
$ws.Range("D45").Value = '17.05'
$ws.Range("E45").Value = '  -1.46%  '
$ws.Range("E46").Value = '  -3.03%  '
$ws.Range("D47").Value = '22.22'
$ws.Range("E47").Value = '  +0.14%  '
$ws.Range("D48").Value = '2.157.51'
$ws.Range("E48").Value = '  +3.56%  '
$ws.Range("E49").Value = '  +1.71%  '
$ws.Range("D50").Value = '2.02'
$ws.Range("E50").Value = '  +21.64%  '
$ws.Range("E51").Value = '  -2.72%  '
